{"js": "// Load all body paragraphs so we can locate the anchor paragraphs by text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the \"foreach (...)\" bullet -- the new sub-bullets are inserted\n// right after it (and right before the \"More elegant way...\" bullet).\nlet foreachPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Access value using temporary\") !== -1 && t.indexOf(\"foreach\") !== -1) {\n    foreachPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!foreachPara) {\n  throw new Error(\"Could not locate the 'foreach (...)' paragraph.\");\n}\n\n// Remove the old (now stale) \"_GoBack\" bookmark -- Word relocates this\n// internal \"last edit position\" bookmark whenever new content is typed.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Insert the three new sub-bullets right after the \"foreach (...)\" bullet.\n// insertParagraph() clones the source paragraph's style/numbering (ilvl 1,\n// numId 8, Listeavsnitt), matching the diff.\nconst p1 = foreachPara.insertParagraph(\n  \"\\u201cArray.Reverse(anArray)\\u201d :  \",\n  Word.InsertLocation.after\n);\nconst p2 = p1.insertParagraph(\n  \"Is a .Net framework method for the class \\u201cArray\\u201d. Takes an array as argument.\",\n  Word.InsertLocation.after\n);\nconst p3 = p2.insertParagraph(\n  \"Part of System Nnamespace.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Re-fetch the freshly-inserted \"Part of System Nnamespace.\" paragraph by\n// text so its Range reflects the just-synced document (a Range obtained\n// from a paragraph object created earlier in the same batch can otherwise\n// point at the wrong offset).\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nlet p3Fresh = null;\nfor (let i = 0; i < refreshedParagraphs.items.length; i++) {\n  if (refreshedParagraphs.items[i].text === \"Part of System Nnamespace.\") {\n    p3Fresh = refreshedParagraphs.items[i];\n    break;\n  }\n}\nif (!p3Fresh) {\n  throw new Error(\"Could not locate the 'Part of System Nnamespace.' paragraph.\");\n}\n\n// Re-plant the \"_GoBack\" bookmark at the end of the new last sub-bullet\n// (right after \"Part of System Nnamespace.\"), mirroring where Word leaves\n// the cursor after typing the new text.\nconst endOfP3 = p3Fresh.getRange(\"End\");\nendOfP3.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"foreach (string name in names) : ...\" bullet -- the new\n# sub-bullets are inserted right after it (and right before the\n# \"More elegant way...\" bullet).\n$findRng = $d.Content\n$findRng.Find.Execute(\"foreach (string name in names\") | Out-Null\n$foreachPara = $findRng.Paragraphs(1)\n\n$foreachIdx = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs($i).Range.Start -eq $foreachPara.Range.Start) {\n    $foreachIdx = $i\n    break\n  }\n}\nif ($foreachIdx -eq 0) {\n  throw \"Could not locate the 'foreach (...)' paragraph.\"\n}\n\n$quote = [char]0x201C\n$unquote = [char]0x201D\n\n# Insert the three new sub-bullets right after the \"foreach (...)\" bullet.\n# InsertParagraphAfter() clones the source paragraph's style/numbering\n# (ilvl 1, numId 8, Listeavsnitt), matching the diff. Re-fetching each\n# paragraph by its (now fixed) index -- rather than reusing a cached Range\n# object -- keeps each Range anchored to the right place as the document\n# grows.\n$d.Paragraphs($foreachIdx).Range.InsertParagraphAfter() | Out-Null\n$d.Paragraphs($foreachIdx + 1).Range.Text = $quote + \"Array.Reverse(anArray)\" + $unquote + \" :  \"\n\n$d.Paragraphs($foreachIdx + 1).Range.InsertParagraphAfter() | Out-Null\n$d.Paragraphs($foreachIdx + 2).Range.Text = \"Is a .Net framework method for the class \" + $quote + \"Array\" + $unquote + \". Takes an array as argument.\"\n\n$d.Paragraphs($foreachIdx + 2).Range.InsertParagraphAfter() | Out-Null\n$d.Paragraphs($foreachIdx + 3).Range.Text = \"Part of System Nnamespace.\"\n\n# Remove the old (now stale) \"_GoBack\" bookmark -- Word relocates this\n# internal \"last edit position\" bookmark whenever new content is typed --\n# and re-plant it at the end of the new last sub-bullet\n# (\"Part of System Nnamespace.\"), mirroring where Word leaves the cursor\n# after typing the new text.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$bookmarkRng = $d.Paragraphs($foreachIdx + 3).Range\n$bookmarkRng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRng) | Out-Null\n"}
